$wb = $excel.ActiveWorkbook

# --- Sheet "y": append new monthly rows (50-61) ---
$wsY = $wb.Worksheets.Item("y")
$wsY.Cells.Item(50, 1).Value = 44256
$wsY.Cells.Item(50, 2).Value = 7001
$wsY.Cells.Item(51, 1).Value = 44287
$wsY.Cells.Item(51, 2).Value = 7332
$wsY.Cells.Item(52, 1).Value = 44317
$wsY.Cells.Item(52, 2).Value = 8201
$wsY.Cells.Item(53, 1).Value = 44348
$wsY.Cells.Item(53, 2).Value = 6932
$wsY.Cells.Item(54, 1).Value = 44378
$wsY.Cells.Item(54, 2).Value = 5988
$wsY.Cells.Item(55, 1).Value = 44409
$wsY.Cells.Item(55, 2).Value = 6450
$wsY.Cells.Item(56, 1).Value = 44440
$wsY.Cells.Item(56, 2).Value = 7923
$wsY.Cells.Item(57, 1).Value = 44470
$wsY.Cells.Item(57, 2).Value = 6021
$wsY.Cells.Item(58, 1).Value = 44501
$wsY.Cells.Item(58, 2).Value = 5302
$wsY.Cells.Item(59, 1).Value = 44531
$wsY.Cells.Item(59, 2).Value = 5542
$wsY.Cells.Item(60, 1).Value = 44562
$wsY.Cells.Item(60, 2).Value = 6450
$wsY.Cells.Item(61, 1).Value = 44593
$wsY.Cells.Item(61, 2).Value = 7203

# --- Sheet "X": append new monthly rows (62-73) ---
$wsX = $wb.Worksheets.Item("X")
$wsX.Cells.Item(62, 1).Value = 44621
$wsX.Cells.Item(62, 2).Value = 3000
$wsX.Cells.Item(62, 3).Value = 6893
$wsX.Cells.Item(63, 1).Value = 44652
$wsX.Cells.Item(63, 2).Value = 3121
$wsX.Cells.Item(63, 3).Value = 6751
$wsX.Cells.Item(64, 1).Value = 44682
$wsX.Cells.Item(64, 2).Value = 3240
$wsX.Cells.Item(64, 3).Value = 6678
$wsX.Cells.Item(65, 1).Value = 44713
$wsX.Cells.Item(65, 2).Value = 2425
$wsX.Cells.Item(65, 3).Value = 6816
$wsX.Cells.Item(66, 1).Value = 44743
$wsX.Cells.Item(66, 2).Value = 2555
$wsX.Cells.Item(66, 3).Value = 6724
$wsX.Cells.Item(67, 1).Value = 44774
$wsX.Cells.Item(67, 2).Value = 2345
$wsX.Cells.Item(67, 3).Value = 6677
$wsX.Cells.Item(68, 1).Value = 44805
$wsX.Cells.Item(68, 2).Value = 2745
$wsX.Cells.Item(68, 3).Value = 6344
$wsX.Cells.Item(69, 1).Value = 44835
$wsX.Cells.Item(69, 2).Value = 2464
$wsX.Cells.Item(69, 3).Value = 6433
$wsX.Cells.Item(70, 1).Value = 44866
$wsX.Cells.Item(70, 2).Value = 2334
$wsX.Cells.Item(70, 3).Value = 6441
$wsX.Cells.Item(71, 1).Value = 44896
$wsX.Cells.Item(71, 2).Value = 2131
$wsX.Cells.Item(71, 3).Value = 6454
$wsX.Cells.Item(72, 1).Value = 44927
$wsX.Cells.Item(72, 2).Value = 1998
$wsX.Cells.Item(72, 3).Value = 6321
$wsX.Cells.Item(73, 1).Value = 44958
$wsX.Cells.Item(73, 2).Value = 2034
$wsX.Cells.Item(73, 3).Value = 6212

# --- Drop the "&L&"CorpoS"&10&K000000Internal&1#" internal-marking header/footer on every sheet ---
$wsInfo = $wb.Worksheets.Item("Info")
foreach ($sheet in @($wsInfo, $wsY, $wsX)) {
    $sheet.PageSetup.LeftHeader = ""
    $sheet.PageSetup.CenterHeader = ""
    $sheet.PageSetup.RightHeader = ""
    $sheet.PageSetup.LeftFooter = ""
    $sheet.PageSetup.CenterFooter = ""
    $sheet.PageSetup.RightFooter = ""
}

# --- Restore the on-disk cell selection / scroll position for each sheet ---
$wsInfo.Activate()
$wsInfo.Range("I15").Select()

$wsY.Activate()
$wsY.Range("F47").Select()

$wsX.Activate()
$wsX.Range("G36").Select()

# "Info" was the active sheet/tab in the source workbook
$wsInfo.Activate()
